$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural edit -------------------------------------------------
# The sheet is a weekly "MarketBeatRank" watch report: column A holds the
# analyst/firm name, and each subsequent column is a weekly snapshot
# (header in row 1 is "Jun_xx"), newest weeks inserted just after column A.
# This revision adds two new weekly snapshot columns ("Jun_17", "Jun_15")
# in front of the existing "Jun_13" / "Jun_10" columns, pushing the old
# data two columns to the right (old B -> D, old C -> E).
$ws.Columns("B:C").Insert()

# --- New column headers ----------------------------------------------
$ws.Cells.Item(1, 2).Value = "Jun_17"
$ws.Cells.Item(1, 3).Value = "Jun_15"

# --- Fill the two new snapshot columns with the default "UN" rating --
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# --- Row 18 (Royal Bank of Canada) already had a highlighted rating
# change note in the old B column (now column D). The newest snapshot
# (column C, "Jun_15") picks up that same rating-change note, with its
# own highlighted fill so it stands out like the other flagged cells.
$ws.Cells.Item(18, 3).Value = $ws.Cells.Item(18, 4).Text
$ws.Cells.Item(18, 3).Interior.ColorIndex = $ws.Cells.Item(18, 4).Interior.ColorIndex

# --- Column widths: keep the original 8-wide formatting, now spread
# across the two newly inserted columns plus the original one.
$ws.Columns("C:E").ColumnWidth = 7.1666666666667
